$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (prices in column D, hourly volume % in column E)
# Values are written as text: NumberFormat "@" prevents Excel from auto-coercing
# numeric-looking strings (e.g. "0.630", "242.80") into floats and dropping
# trailing zeros; resetting the Style back to "Normal" afterwards avoids leaving
# any stray per-cell formatting behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.392.56'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.07%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.800.63'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.03%  '
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '225.26'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.604'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.36%  '
$ws.Range("E7").Value = '  -0.20%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '36.03'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.22%  '
$ws.Range("E9").Value = '  -2.32%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0678'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0965'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.41%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.060.86'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.02%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.22'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.54%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.810.43'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.630'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.85%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '34.374.46'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("E17").Value = '  +2.39%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.51'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.98%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '242.80'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.81%  '
$ws.Range("E20").Value = '  -2.79%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.26'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.06%  '
$ws.Range("E22").Value = '  -0.29%  '
$ws.Range("E23").Value = '  -1.51%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.20'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.68%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '170.02'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.53%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.87'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.34'
$ws.Range("D27").Style = "Normal"
$ws.Range("E28").Value = '  +2.45%  '
$ws.Range("E29").Value = '  -0.31%  '
$ws.Range("E30").Value = '  -1.70%  '
$ws.Range("E31").Value = '  -0.38%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.22'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.54%  '
$ws.Range("E33").Value = '  -2.14%  '
$ws.Range("E34").Value = '  -3.09%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.363.86'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.56%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.649'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.95%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.05'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.49%  '
$ws.Range("E38").Value = '  -6.88%  '
$ws.Range("E39").Value = '  -1.41%  '
$ws.Range("E40").Value = '  +0.84%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '81.09'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.17%  '
$ws.Range("E42").Value = '  -1.90%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.936'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.19%  '
$ws.Range("E44").Value = '  +5.09%  '
$ws.Range("E45").Value = '  -2.90%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0498'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.70%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.962.95'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.08%  '
$ws.Range("E48").Value = '  -3.21%  '
$ws.Range("E49").Value = '  -0.29%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '101.94'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.37%  '
$ws.Range("E51").Value = '  -4.17%  '
